$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'315.41"
$ws.Range("E2").Value = "'3.26%"
$ws.Range("D3").Value = "'35.63"
$ws.Range("E3").Value = "'-1.44%"
$ws.Range("D4").Value = "'5.120"
$ws.Range("E4").Value = "'0.43%"
$ws.Range("D5").Value = "'0.08110"
$ws.Range("E5").Value = "'2.94%"
$ws.Range("E6").Value = "'-0.26%"
$ws.Range("D7").Value = "'8.003"
$ws.Range("E7").Value = "'0.99%"
$ws.Range("D8").Value = "'4.150"
$ws.Range("E8").Value = "'1.05%"
$ws.Range("D9").Value = "'0.9256"
$ws.Range("E9").Value = "'0.36%"
$ws.Range("D10").Value = "'0.1026"
$ws.Range("E10").Value = "'5.49%"
$ws.Range("D11").Value = "'0.1876"
$ws.Range("E11").Value = "'1.06%"
$ws.Range("D12").Value = "'0.09199"
$ws.Range("E12").Value = "'6.60%"
$ws.Range("D13").Value = "'0.03589"
$ws.Range("E13").Value = "'0.84%"
$ws.Range("D14").Value = "'0.09909"
$ws.Range("E14").Value = "'-0.15%"
$ws.Range("D15").Value = "'0.001431"
$ws.Range("D16").Value = "'0.005679"
$ws.Range("E16").Value = "'1.01%"
$ws.Range("D17").Value = "'3.477"
$ws.Range("E17").Value = "'0.36%"
$ws.Range("D18").Value = "'2.828"
$ws.Range("E18").Value = "'7.29%"
$ws.Range("E19").Value = "'-0.87%"
$ws.Range("D20").Value = "'0.1331"
$ws.Range("E20").Value = "'1.04%"
$ws.Range("D21").Value = "'5.129"
$ws.Range("E21").Value = "'-0.59%"
$ws.Range("D22").Value = "'0.2219"
$ws.Range("D23").Value = "'0.04573"
$ws.Range("E23").Value = "'0.54%"
$ws.Range("D24").Value = "'0.001246"
$ws.Range("E24").Value = "'1.06%"
$ws.Range("D25").Value = "'0.004708"
$ws.Range("E25").Value = "'-6.85%"
$ws.Range("D26").Value = "'0.0001251"
$ws.Range("E26").Value = "'-21.94%"
$ws.Range("D27").Value = "'0.0004504"
$ws.Range("E27").Value = "'-5.04%"
$ws.Range("D39").Value = "'0.01963"
$ws.Range("E39").Value = "'6.31%"
$ws.Range("D40").Value = "'0.04862"
$ws.Range("E40").Value = "'2.24%"
$ws.Range("D41").Value = "'0.007716"
$ws.Range("E41").Value = "'1.93%"
$ws.Range("E42").Value = "'-0.61%"
$ws.Range("D43").Value = "'0.007828"
$ws.Range("E43").Value = "'1.39%"
$ws.Range("D44").Value = "'0.002152"
$ws.Range("E44").Value = "'-2.94%"
$ws.Range("D45").Value = "'0.01164"
$ws.Range("E45").Value = "'2.77%"
$ws.Range("D46").Value = "'0.00006516"
$ws.Range("E46").Value = "'3.11%"
$ws.Range("E47").Value = "'0.22%"
$ws.Range("D48").Value = "'36.65"
$ws.Range("E48").Value = "'-22.78%"
$ws.Range("D49").Value = "'0.001701"
$ws.Range("E49").Value = "'-14.80%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.22%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.22%"
